$wb = $excel.ActiveWorkbook

$wsSummary = $wb.Worksheets.Item("summary")
$wsModelFit = $wb.Worksheets.Item("model_fit")
$wsSteps = $wb.Worksheets.Item("steps")

# summary
$wsSummary.Range("C2").Value = 714
$wsSummary.Range("D2").Value = 689
$wsSummary.Range("E2").Value = 80.99
$wsSummary.Range("H2").Value = 1.04
$wsSummary.Range("I2").Value = 0.66
$wsSummary.Range("K2").Value = 0.07
$wsSummary.Range("L2").Value = 0.86
$wsSummary.Range("C3").Value = 714
$wsSummary.Range("D3").Value = 681
$wsSummary.Range("F3").Value = -1.49
$wsSummary.Range("H3").Value = 0.92
$wsSummary.Range("I3").Value = -1.02
$wsSummary.Range("J3").Value = 0.33
$wsSummary.Range("K3").Value = 0.06
$wsSummary.Range("L3").Value = 1.01
$wsSummary.Range("C4").Value = 714
$wsSummary.Range("D4").Value = 689
$wsSummary.Range("E4").Value = 48.62
$wsSummary.Range("F4").Value = 0.07
$wsSummary.Range("H4").Value = 1.04
$wsSummary.Range("I4").Value = 1.03
$wsSummary.Range("J4").Value = 0.28
$wsSummary.Range("K4").Value = 0.05
$wsSummary.Range("L4").Value = 0.97
$wsSummary.Range("C5").Value = 714
$wsSummary.Range("D5").Value = 687
$wsSummary.Range("F5").Value = 0.08
$wsSummary.Range("I5").Value = -2.44
$wsSummary.Range("K5").Value = 0.08
$wsSummary.Range("C6").Value = 714
$wsSummary.Range("D6").Value = 669
$wsSummary.Range("E6").Value = 18.54
$wsSummary.Range("F6").Value = 1.85
$wsSummary.Range("I6").Value = 1.08
$wsSummary.Range("J6").Value = 0.2
$wsSummary.Range("L6").Value = 0.74
$wsSummary.Range("C7").Value = 714
$wsSummary.Range("D7").Value = 676
$wsSummary.Range("E7").Value = 63.76
$wsSummary.Range("F7").Value = -0.71
$wsSummary.Range("H7").Value = 1
$wsSummary.Range("I7").Value = -0.06
$wsSummary.Range("J7").Value = 0.34
$wsSummary.Range("K7").Value = 0.06
$wsSummary.Range("L7").Value = 1.14
$wsSummary.Range("C8").Value = 714
$wsSummary.Range("D8").Value = 643
$wsSummary.Range("E8").Value = 55.37
$wsSummary.Range("F8").Value = -0.27
$wsSummary.Range("I8").Value = 0.09
$wsSummary.Range("J8").Value = 0.38
$wsSummary.Range("L8").Value = 1.32
$wsSummary.Range("C9").Value = 714
$wsSummary.Range("D9").Value = 623
$wsSummary.Range("E9").Value = 51.04
$wsSummary.Range("F9").Value = -0.06
$wsSummary.Range("I9").Value = 0.27
$wsSummary.Range("J9").Value = 0.32
$wsSummary.Range("L9").Value = 1.03
$wsSummary.Range("C10").Value = 714
$wsSummary.Range("E10").Value = 47.62
$wsSummary.Range("F10").Value = 0.13
$wsSummary.Range("G10").Value = 0.1
$wsSummary.Range("H10").Value = 1.03
$wsSummary.Range("I10").Value = 0.78
$wsSummary.Range("J10").Value = 0.32
$wsSummary.Range("L10").Value = 1
$wsSummary.Range("C11").Value = 714
$wsSummary.Range("E11").Value = 42.47
$wsSummary.Range("F11").Value = 0.37
$wsSummary.Range("H11").Value = 1.03
$wsSummary.Range("I11").Value = 0.8
$wsSummary.Range("J11").Value = 0.33
$wsSummary.Range("K11").Value = 0.03
$wsSummary.Range("L11").Value = 0.99
$wsSummary.Range("C12").Value = 714
$wsSummary.Range("D12").Value = 336
$wsSummary.Range("F12").Value = -0.25
$wsSummary.Range("H12").Value = 0.92
$wsSummary.Range("I12").Value = -1.26
$wsSummary.Range("J12").Value = 0.56
$wsSummary.Range("K12").Value = 0.1
$wsSummary.Range("L12").Value = 0.9

# model_fit
$wsModelFit.Range("B2").Value = 714
$wsModelFit.Range("D2").Value = 8871
$wsModelFit.Range("E2").Value = 8903
$wsModelFit.Range("F2").Value = 8976
$wsModelFit.Range("G2").Value = 0.675
$wsModelFit.Range("H2").Value = 0.581
$wsModelFit.Range("B3").Value = 714
$wsModelFit.Range("D3").Value = 8823
$wsModelFit.Range("E3").Value = 8875
$wsModelFit.Range("F3").Value = 8994
$wsModelFit.Range("G3").Value = 0.692
$wsModelFit.Range("H3").Value = 0.583

# steps
$wsSteps.Range("B2").Value = "0.17 (0.112)"
$wsSteps.Range("C2").Value = "'-0.17"
$wsSteps.Range("C2").Style = "Normal"
$wsSteps.Range("B4").Value = "1.39 (0.122)"
$wsSteps.Range("C4").Value = "-1.03 (0.142)"
$wsSteps.Range("D4").Value = "'-0.36"
$wsSteps.Range("D4").Style = "Normal"
